# Add two new columns, I ("I0") and J ("IF"), to the right of the existing
# H ("IP") column on the active sheet, matching the header style used by
# the other header cells and filling in the per-row numeric data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---------------------------------------------------
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the header formatting (bold font, centered, thin border) used by
# the existing header cells (e.g. H1) by copying its format onto I1:J1.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# --- Data rows (rows 2-21) -------------------------------------------------
$iValues = @(6, 8, 8, 4, 8, 8, 4, 6, 4, 7, 6, 1, 2, 4, 9, 3, 4, 5, 4, 8)
$jValues = @(8, 8, 8, 6, 9, 8, 4, 7, 5, 7, 8, 4, 4, 5, 9, 5, 5, 6, 4, 9)

for ($r = 0; $r -lt $iValues.Length; $r++) {
    $row = $r + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$r]
    $ws.Cells.Item($row, 10).Value = $jValues[$r]
}
